$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3969039.5
$ws.Range("I18").Value = 3969039.5
$ws.Range("K18").Value = 3969039.5
$ws.Range("M18").Value = -3968755.5
$ws.Range("H19").Value = 6803209.5
$ws.Range("I19").Value = 14286179
$ws.Range("K19").Value = 14286179
$ws.Range("M19").Value = -14286004
$ws.Range("H33").Value = 506.4643
$ws.Range("I33").Value = 391.6842
$ws.Range("J33").Value = 748.7778
$ws.Range("K33").Value = 391.6842
$ws.Range("L33").Value = 748.7778
$ws.Range("M33").Value = -162.6842
$ws.Range("N33").Value = -1206.7778
$ws.Range("H129").Value = 4099534
$ws.Range("J129").Value = 1187.0566
$ws.Range("L129").Value = 3561.1698
$ws.Range("N129").Value = -13561.1698
$ws.Range("H135").Value = 719.25
$ws.Range("I135").Value = 513.63635
$ws.Range("J135").Value = 1473.1666
$ws.Range("K135").Value = 4622.72715
$ws.Range("L135").Value = 13258.4994
$ws.Range("M135").Value = -2087.72715
$ws.Range("N135").Value = -18328.4994
$ws.Range("H139").Value = 27785.908
$ws.Range("J139").Value = 29330
$ws.Range("L139").Value = 29330
$ws.Range("N139").Value = -39610
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8335130
$ws.Range("I2").Value = 11364830
$ws.Range("J2").Value = 3453.25
$ws.Range("K2").Value = 11364830
$ws.Range("L2").Value = 3453.25
$ws.Range("M2").Value = -11364717
$ws.Range("N2").Value = -3679.25
$ws.Range("H45").Value = 1330.6316
$ws.Range("I45").Value = 1063.6364
$ws.Range("J45").Value = 3092.8
$ws.Range("K45").Value = 1063.6364
$ws.Range("L45").Value = 3092.8
$ws.Range("M45").Value = -686.6364000000001
$ws.Range("N45").Value = -3846.8
$ws.Range("H97").Value = 565.4167
$ws.Range("I97").Value = 536.38464
$ws.Range("J97").Value = 640.9
$ws.Range("K97").Value = 536.38464
$ws.Range("L97").Value = 640.9
$ws.Range("M97").Value = -40.38463999999999
$ws.Range("N97").Value = -1632.9
$ws.Range("H116").Value = 8335130
$ws.Range("I116").Value = 11364830
$ws.Range("J116").Value = 3453.25
$ws.Range("K116").Value = 11364830
$ws.Range("L116").Value = 3453.25
$ws.Range("M116").Value = -11362536
$ws.Range("N116").Value = -8041.25
$ws.Range("H132").Value = 2024.7377
$ws.Range("I132").Value = 1402.921
$ws.Range("J132").Value = 3052.087
$ws.Range("K132").Value = 4208.763
$ws.Range("L132").Value = 9156.261
$ws.Range("M132").Value = -1678.763
$ws.Range("N132").Value = -14216.261
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8335130
$ws.Range("I3").Value = 11364830
$ws.Range("J3").Value = 3453.25
$ws.Range("K3").Value = 11364830
$ws.Range("L3").Value = 3453.25
$ws.Range("M3").Value = -11364716
$ws.Range("N3").Value = -3681.25
$ws.Range("H94").Value = 627.46155
$ws.Range("I94").Value = 577
$ws.Range("J94").Value = 728.38464
$ws.Range("K94").Value = 577
$ws.Range("L94").Value = 728.38464
$ws.Range("M94").Value = -126
$ws.Range("N94").Value = -1630.38464
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1996.9
$ws.Range("I22").Value = 505
$ws.Range("J22").Value = 2991.5
$ws.Range("K22").Value = 505
$ws.Range("L22").Value = 2991.5
$ws.Range("M22").Value = -155
$ws.Range("N22").Value = -3691.5
$ws.Range("H25").Value = 47008.668
$ws.Range("I25").Value = 1000
$ws.Range("K25").Value = 1000
$ws.Range("M25").Value = -826
$ws.Range("H132").Value = 1479.862
$ws.Range("I132").Value = 1093.2046
$ws.Range("J132").Value = 2695.0715
$ws.Range("K132").Value = 3279.6138
$ws.Range("L132").Value = 8085.2145
$ws.Range("M132").Value = -749.6138000000001
$ws.Range("N132").Value = -13145.2145
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1447.4375
$ws.Range("I5").Value = 858.5454999999999
$ws.Range("J5").Value = 2743
$ws.Range("K5").Value = 2575.6365
$ws.Range("L5").Value = 8229
$ws.Range("M5").Value = -2463.6365
$ws.Range("N5").Value = -8453
$ws.Range("H87").Value = 9875.5
$ws.Range("J87").Value = 16000
$ws.Range("L87").Value = 48000
$ws.Range("N87").Value = -50496
$ws.Range("H90").Value = 9875.5
$ws.Range("J90").Value = 16000
$ws.Range("L90").Value = 144000
$ws.Range("N90").Value = -156480
$ws.Range("H131").Value = 1413
$ws.Range("I131").Value = 1707.5
$ws.Range("J131").Value = 1322.3846
$ws.Range("K131").Value = 5122.5
$ws.Range("L131").Value = 3967.1538
$ws.Range("M131").Value = -82.5
$ws.Range("N131").Value = -14047.1538
$ws.Range("H135").Value = 1447.4375
$ws.Range("I135").Value = 858.5454999999999
$ws.Range("J135").Value = 2743
$ws.Range("K135").Value = 7726.9095
$ws.Range("L135").Value = 24687
$ws.Range("M135").Value = -5191.9095
$ws.Range("N135").Value = -29757
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2482.889
$ws.Range("I126").Value = 1631.6666
$ws.Range("J126").Value = 3334.111
$ws.Range("K126").Value = 4894.9998
$ws.Range("L126").Value = 10002.333
$ws.Range("M126").Value = -2424.9998
$ws.Range("N126").Value = -14942.333
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 66667896
$ws.Range("I22").Value = 166667170
$ws.Range("J22").Value = 1714.7778
$ws.Range("K22").Value = 166667170
$ws.Range("L22").Value = 1714.7778
$ws.Range("M22").Value = -166666875
$ws.Range("N22").Value = -2304.7778
$ws.Range("H27").Value = 66667896
$ws.Range("I27").Value = 166667170
$ws.Range("J27").Value = 1714.7778
$ws.Range("K27").Value = 166667170
$ws.Range("L27").Value = 1714.7778
$ws.Range("M27").Value = -166667063
$ws.Range("N27").Value = -1928.7778
$ws.Range("H55").Value = 854.0526
$ws.Range("I55").Value = 197
$ws.Range("J55").Value = 1445.4
$ws.Range("K55").Value = 197
$ws.Range("L55").Value = 1445.4
$ws.Range("M55").Value = -24
$ws.Range("N55").Value = -1791.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11030.5
$ws.Range("I132").Value = 2420.4883
$ws.Range("J132").Value = 35712.535
$ws.Range("K132").Value = 7261.4649
$ws.Range("L132").Value = 107137.605
$ws.Range("M132").Value = -4731.4649
$ws.Range("N132").Value = -112197.605
